# Workbook: "Hortaliza, Femacal de La Calera - Zapallo"
# Two new weekly price rows are inserted right after row 190 (before the
# existing row that used to be 191), pushing all subsequent rows down by 2.
# New dimension becomes A1:R296 (was A1:R294).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 191:192 - everything below shifts down by 2.
$ws.Rows("191:192").Insert()

# New row 191: Camote, 1a (guarda)
$ws.Cells.Item(191, 1).Value2 = 3
$ws.Cells.Item(191, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(191, 3).Value2 = "Coquimbo"
$ws.Cells.Item(191, 4).Value2 = 44452
$ws.Cells.Item(191, 5).Value2 = 5
$ws.Cells.Item(191, 6).Value2 = 100112045
$ws.Cells.Item(191, 7).Value2 = "Zapallo"
$ws.Cells.Item(191, 8).Value2 = "Camote"
$ws.Cells.Item(191, 9).Value2 = "1a (guarda)"
$ws.Cells.Item(191, 10).Value2 = 310
$ws.Cells.Item(191, 11).Value2 = 750
$ws.Cells.Item(191, 12).Value2 = 800
$ws.Cells.Item(191, 13).Value2 = 771
$ws.Cells.Item(191, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(191, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(191, 16).Value2 = 771
$ws.Cells.Item(191, 17).Value2 = 1
$ws.Cells.Item(191, 18).Value2 = "Hortaliza"

# New row 192: Camote, 2a (guarda)
$ws.Cells.Item(192, 1).Value2 = 3
$ws.Cells.Item(192, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(192, 3).Value2 = "Coquimbo"
$ws.Cells.Item(192, 4).Value2 = 44452
$ws.Cells.Item(192, 5).Value2 = 5
$ws.Cells.Item(192, 6).Value2 = 100112045
$ws.Cells.Item(192, 7).Value2 = "Zapallo"
$ws.Cells.Item(192, 8).Value2 = "Camote"
$ws.Cells.Item(192, 9).Value2 = "2a (guarda)"
$ws.Cells.Item(192, 10).Value2 = 150
$ws.Cells.Item(192, 11).Value2 = 600
$ws.Cells.Item(192, 12).Value2 = 600
$ws.Cells.Item(192, 13).Value2 = 600
$ws.Cells.Item(192, 14).Value2 = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(192, 15).Value2 = "Provincia de Talca"
$ws.Cells.Item(192, 16).Value2 = 600
$ws.Cells.Item(192, 17).Value2 = 1
$ws.Cells.Item(192, 18).Value2 = "Hortaliza"

# Preserve the date-format (same as the rest of column D / "Fecha") used by
# the neighbouring rows for the two new cells.
$dateFmt = $ws.Cells.Item(193, 4).NumberFormat
$ws.Cells.Item(191, 4).NumberFormat = $dateFmt
$ws.Cells.Item(192, 4).NumberFormat = $dateFmt
